$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row contents (columns B, E through AD) between pairs of rows.
# Columns A, C, D keep their values (row index, Div, Date).
$pairs = @(
    @(47, 48),
    @(71, 72),
    @(101, 102),
    @(149, 150),
    @(232, 233),
    @(248, 249),
    @(271, 272),
    @(307, 308)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Column B (col 2)
    $b1 = $ws.Cells.Item($r1, 2).Value2
    $b2 = $ws.Cells.Item($r2, 2).Value2
    $ws.Cells.Item($r1, 2).Value2 = $b2
    $ws.Cells.Item($r2, 2).Value2 = $b1

    # Columns E (5) through AD (30)
    for ($c = 5; $c -le 30; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value2 = $v2
        $ws.Cells.Item($r2, $c).Value2 = $v1
    }
}
